$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 379, shifting existing rows 379:407 down to 380:408
$ws.Rows("379:379").Insert(-4121)

# Populate the newly inserted row 379 with the new weekly record
$ws.Cells.Item(379, 1).Value = 5
$ws.Cells.Item(379, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(379, 3).Value = "Maule"
$ws.Cells.Item(379, 4).Value = 45021
$ws.Cells.Item(379, 5).Value = 7
$ws.Cells.Item(379, 6).Value = 100112009
$ws.Cells.Item(379, 7).Value = "Acelga"
$ws.Cells.Item(379, 8).Value = "Sin especificar"
$ws.Cells.Item(379, 9).Value = "Primera"
$ws.Cells.Item(379, 10).Value = 400
$ws.Cells.Item(379, 11).Value = 2500
$ws.Cells.Item(379, 12).Value = 2500
$ws.Cells.Item(379, 13).Value = 2500
$ws.Cells.Item(379, 14).Value = "$/docena de atados (4 kilos)"
$ws.Cells.Item(379, 15).Value = "Región del Maule"
$ws.Cells.Item(379, 16).Value = 625
$ws.Cells.Item(379, 17).Value = 4
$ws.Cells.Item(379, 18).Value = "Hortaliza"
